$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117 (existing rows 117-131 shift down to 118-132),
# inheriting the formatting (styles) of the row above it.
$ws.Rows("117:117").Insert()

# Preserve the standard row height used throughout this block of shortcuts.
$ws.Rows("117:117").RowHeight = 17

# Populate the newly inserted row with the new "Estimate Noise" (EN) shortcut.
$ws.Range("A117").Value = "Estimate Noise"
$ws.Range("B117").Value = "EN"

# The printed area grew by one row because of the inserted row.
$ws.PageSetup.PrintArea = "`$A`$1:`$C`$123"

# Update the active selection to match the edited location.
$ws.Range("C118").Select()
